# The "ReportingOrganisationGroup" sheet has columns:
#   A=code, B=name, C=status, D=codeforiati:group-code, E=codeforiati:group-name
# This edit swaps the contents of the group-code (D) and group-name (E)
# columns for every row (including the header row), so that column D now
# holds the group name and column E now holds the group code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
